# ntp_template_map.xlsx edit:
#  - "studies" block gets a new "Dose Frequency (per Study)" mapping row,
#    while the original "Dose Frequency (Per Study)" -> dose_frequency mapping
#    is kept as a separate row right below it (handling multiple weight/dose
#    frequency fields).
#  - the duplicate "Initial Body Weight (g)" row in the "subjects" block is
#    removed (de-duping subjects).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after row 29 (the current "Dose Frequency (Per
# Study)" row in the studies block). This shifts the old row 29 -> row 30,
# row 30 -> row 31, etc.
$ws.Rows("30:30").Insert()

# Duplicate the original row 29 content (A/B/C) down into the newly inserted
# row 30, so the "Dose Frequency (Per Study)" -> dose_frequency mapping is
# preserved as its own row.
$ws.Range("A29:C29").Copy()
$ws.Range("A30").PasteSpecial()
$excel.CutCopyMode = $false

# Row 29 now becomes the new, second dose-frequency column mapping.
$ws.Range("B29").Value = "Dose Frequency (per Study)"

# After the insert above, the old "Initial Body Weight (g)" row (previously
# row 56 in the subjects block) has shifted down to row 57. Remove it so the
# subjects block no longer has the duplicate weight field.
$ws.Rows("57:57").Delete()

# Restore the view state (selection / scroll position) recorded in the
# edited workbook.
$ws.Range("B56").Select()
$excel.ActiveWindow.ScrollRow = 39
